$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "time_taken" in F1, matching the style of the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Values for time_taken column, rows 2..14
$times = @(
    "2021-10-05 13:41:37.059123",
    "2021-10-05 13:41:37.059135",
    "2021-10-05 13:41:37.059139",
    "2021-10-05 13:41:37.059142",
    "2021-10-05 13:41:37.059145",
    "2021-10-05 13:41:37.059149",
    "2021-10-05 13:41:37.059152",
    "2021-10-05 13:41:37.059155",
    "2021-10-05 13:41:37.059158",
    "2021-10-05 13:41:37.059161",
    "2021-10-05 13:41:37.059164",
    "2021-10-05 13:41:37.059167",
    "2021-10-05 13:41:37.059170"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
